$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19 (id 7015)
$ws.Range("H19").Value = 896.64514
$ws.Range("J19").Value = 1131.3158
$ws.Range("L19").Value = 1131.3158
$ws.Range("N19").Value = -1481.3158
# row 28 (id 27772)
$ws.Range("H28").Value = 2066.4644
$ws.Range("I28").Value = 601.6923
$ws.Range("K28").Value = 601.6923
$ws.Range("M28").Value = -116.6923
# row 34 (id 2160)
$ws.Range("H34").Value = 2822.6365
$ws.Range("I34").Value = 2822.6365
$ws.Range("K34").Value = 2822.6365
$ws.Range("M34").Value = -2619.6365
# row 36 (id 2160)
$ws.Range("H36").Value = 2822.6365
$ws.Range("I36").Value = 2822.6365
$ws.Range("K36").Value = 2822.6365
$ws.Range("M36").Value = -2107.6365
# row 40 (id 5505)
$ws.Range("H40").Value = 4634.3423
$ws.Range("I40").Value = 2678.8572
$ws.Range("J40").Value = 5075.9033
$ws.Range("K40").Value = 2678.8572
$ws.Range("L40").Value = 5075.9033
$ws.Range("M40").Value = -2503.8572
$ws.Range("N40").Value = -5425.9033
# row 54 (id 2174)
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").Value = $null
# row 64 (id 5506)
$ws.Range("H64").Value = 8887.286
$ws.Range("J64").Value = 6396.5
$ws.Range("L64").Value = 6396.5
$ws.Range("N64").Value = -6892.5
# row 67 (id 5506)
$ws.Range("H67").Value = 8887.286
$ws.Range("J67").Value = 6396.5
$ws.Range("L67").Value = 6396.5
$ws.Range("N67").Value = -8112.5
# row 69 (id 12616)
$ws.Range("H69").Value = 500004000
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = $null
# row 72 (id 12616)
$ws.Range("H72").Value = 500004000
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = $null
# row 74 (id 5507)
$ws.Range("H74").Value = 4817.3335
$ws.Range("I74").Value = 4817.3335
$ws.Range("K74").Value = 4817.3335
$ws.Range("M74").Value = -3881.3335
# row 77 (id 5507)
$ws.Range("H77").Value = 4817.3335
$ws.Range("I77").Value = 4817.3335
$ws.Range("K77").Value = 24086.6675
$ws.Range("M77").Value = -19406.6675
# row 80 (id 12605)
$ws.Range("H80").Value = 418.8
$ws.Range("I80").Value = 254.55556
$ws.Range("K80").Value = 763.66668
$ws.Range("M80").Value = 234.33332
# row 83 (id 12605)
$ws.Range("H83").Value = 418.8
$ws.Range("I83").Value = 254.55556
$ws.Range("K83").Value = 2291.00004
$ws.Range("M83").Value = 2700.99996
# row 98 (id 36237)
$ws.Range("H98").Value = 1049.8667
$ws.Range("I98").Value = 1097.3658
$ws.Range("K98").Value = 1097.3658
$ws.Range("M98").Value = 400.6342
# row 100 (id 19906)
$ws.Range("H100").Value = 1074.8889
$ws.Range("I100").Value = 785.17645
$ws.Range("K100").Value = 785.17645
$ws.Range("M100").Value = -244.17645
# row 113 (id 27775)
$ws.Range("H113").Value = 7472.5
$ws.Range("J113").Value = 7755
$ws.Range("L113").Value = 7755
$ws.Range("N113").Value = -14263
# row 122 (id 36237)
$ws.Range("H122").Value = 1049.8667
$ws.Range("I122").Value = 1097.3658
$ws.Range("K122").Value = 3292.0974
$ws.Range("M122").Value = -842.0974000000001
# row 132 (id 44049)
$ws.Range("H132").Value = 45458180
$ws.Range("I132").Value = 66669932
$ws.Range("K132").Value = 200009796
$ws.Range("M132").Value = -200007266
# row 135 (id 44047)
$ws.Range("H135").Value = 1048.8928
$ws.Range("J135").Value = 1966.1111
$ws.Range("L135").Value = 17694.9999
$ws.Range("N135").Value = -22764.9999
# row 137 (id 44013)
$ws.Range("H137").Value = 199435.89
$ws.Range("I137").Value = 445440.5
$ws.Range("J137").Value = 2632.2
$ws.Range("K137").Value = 1336321.5
$ws.Range("L137").Value = 7896.599999999999
$ws.Range("M137").Value = -1333771.5
$ws.Range("N137").Value = -12996.6
# row 138 (id 44169)
$ws.Range("H138").Value = 4015.6365
$ws.Range("I138").Value = 3657.611
$ws.Range("J138").Value = 5626.75
$ws.Range("K138").Value = 10972.833
$ws.Range("L138").Value = 16880.25
$ws.Range("M138").Value = -5832.832999999999
$ws.Range("N138").Value = -27160.25
# row 141 (id 44161)
$ws.Range("H141").Value = 24714
$ws.Range("I141").Value = 11128.2
$ws.Range("J141").Value = 70000
$ws.Range("K141").Value = 33384.60000000001
$ws.Range("L141").Value = 210000
$ws.Range("M141").Value = -28204.60000000001
$ws.Range("N141").Value = -220360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32 (id 44147)
$ws.Range("H32").Value = 2663.238
$ws.Range("I32").Value = 2299.8928
$ws.Range("K32").Value = 2299.8928
$ws.Range("M32").Value = -2012.8928
# row 45 (id 27714)
$ws.Range("H45").Value = 4798126.5
$ws.Range("I45").Value = 6851448
$ws.Range("K45").Value = 6851448
$ws.Range("M45").Value = -6851071
# row 57 (id 39767)
$ws.Range("H57").Value = 6241.5
$ws.Range("I57").Value = 6241.5
$ws.Range("K57").Value = 6241.5
$ws.Range("M57").Value = -5757.5
# row 61 (id 43999)
$ws.Range("H61").Value = 3923.2964
$ws.Range("I61").Value = 3882.077
$ws.Range("K61").Value = 3882.077
$ws.Range("M61").Value = -3670.077
# row 74 (id 44000)
$ws.Range("H74").Value = 59432.67
$ws.Range("I74").Value = 5124.385
$ws.Range("K74").Value = 5124.385
$ws.Range("M74").Value = -4250.385
# row 77 (id 44000)
$ws.Range("H77").Value = 59432.67
$ws.Range("I77").Value = 5124.385
$ws.Range("K77").Value = 25621.925
$ws.Range("M77").Value = -21253.925
# row 88 (id 12530)
$ws.Range("H88").Value = 610
$ws.Range("I88").Value = 610
$ws.Range("K88").Value = 610
$ws.Range("M88").Value = -204
# row 91 (id 12530)
$ws.Range("H91").Value = 610
$ws.Range("I91").Value = 610
$ws.Range("K91").Value = 610
$ws.Range("M91").Value = 794
# row 97 (id 19941)
$ws.Range("H97").Value = 2022886.8
$ws.Range("I97").Value = 2157685.8
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 2157685.8
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -2157189.8
$ws.Range("N97").Value = -1892
# row 102 (id 19945)
$ws.Range("H102").Value = 6414342
$ws.Range("I102").Value = 6947953.5
$ws.Range("K102").Value = 6947953.5
$ws.Range("M102").Value = -6946331.5
# row 107 (id 25645)
$ws.Range("H107").Value = 34400
$ws.Range("J107").Value = 34400
$ws.Range("L107").Value = 34400
$ws.Range("N107").Value = -42080
# row 122 (id 36168)
$ws.Range("H122").Value = 1229450.5
$ws.Range("I122").Value = 3682.375
$ws.Range("J122").Value = 2319022
$ws.Range("K122").Value = 11047.125
$ws.Range("L122").Value = 6957066
$ws.Range("M122").Value = -8597.125
$ws.Range("N122").Value = -6961966
# row 126 (id 39766)
$ws.Range("H126").Value = 6250
$ws.Range("I126").Value = 6250
$ws.Range("K126").Value = 18750
$ws.Range("M126").Value = -16280
# row 132 (id 43997)
$ws.Range("H132").Value = 2563.6216
$ws.Range("I132").Value = 1836.5938
$ws.Range("J132").Value = 7216.6
$ws.Range("K132").Value = 5509.7814
$ws.Range("L132").Value = 21649.8
$ws.Range("M132").Value = -2979.7814
$ws.Range("N132").Value = -26709.8
# row 136 (id 43999)
$ws.Range("H136").Value = 3923.2964
$ws.Range("I136").Value = 3882.077
$ws.Range("K136").Value = 11646.231
$ws.Range("M136").Value = -9096.231

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20 (id 14149)
$ws.Range("H20").Value = 15875268
$ws.Range("I20").Value = 25642446
$ws.Range("J20").Value = 3604.25
$ws.Range("K20").Value = 25642446
$ws.Range("L20").Value = 3604.25
$ws.Range("M20").Value = -25642199
$ws.Range("N20").Value = -4098.25
# row 64 (id 14184)
$ws.Range("H64").Value = 415
$ws.Range("J64").Value = 396.66666
$ws.Range("L64").Value = 396.66666
$ws.Range("N64").Value = -846.66666
# row 67 (id 14184)
$ws.Range("H67").Value = 415
$ws.Range("J67").Value = 396.66666
$ws.Range("L67").Value = 396.66666
$ws.Range("N67").Value = -1956.66666
# row 86 (id 12526)
$ws.Range("H86").Value = 5011727.5
$ws.Range("J86").Value = 3957.3333
$ws.Range("L86").Value = 3957.3333
$ws.Range("N86").Value = -6203.3333
# row 89 (id 12526)
$ws.Range("H89").Value = 5011727.5
$ws.Range("J89").Value = 3957.3333
$ws.Range("L89").Value = 19786.6665
$ws.Range("N89").Value = -31018.6665
# row 94 (id 19939)
$ws.Range("H94").Value = 4333342
$ws.Range("I94").Value = 5347828
$ws.Range("J94").Value = 21776.25
$ws.Range("K94").Value = 5347828
$ws.Range("L94").Value = 21776.25
$ws.Range("M94").Value = -5347377
$ws.Range("N94").Value = -22678.25
# row 99 (id 19943)
$ws.Range("H99").Value = 4929721
$ws.Range("I99").Value = 6806008
$ws.Range("K99").Value = 6806008
$ws.Range("M99").Value = -6804510
# row 105 (id 19947)
$ws.Range("H105").Value = 3736349
$ws.Range("I105").Value = 4234343
$ws.Range("K105").Value = 4234343
$ws.Range("M105").Value = -4232596
# row 107 (id 27706)
$ws.Range("H107").Value = 3403938.5
$ws.Range("I107").Value = 4763766.5
$ws.Range("K107").Value = 4763766.5
$ws.Range("M107").Value = -4761846.5
# row 113 (id 39768)
$ws.Range("H113").Value = 1447922.8
$ws.Range("I113").Value = 1447922.8
$ws.Range("K113").Value = 1447922.8
$ws.Range("M113").Value = -1445752.8
# row 134 (id 43998)
$ws.Range("H134").Value = 2531.532
$ws.Range("I134").Value = 1254.35
$ws.Range("J134").Value = 9829.714
$ws.Range("K134").Value = 3763.05
$ws.Range("L134").Value = 29489.142
$ws.Range("M134").Value = -1228.05
$ws.Range("N134").Value = -34559.142

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16 (id 27691)
$ws.Range("H16").Value = 3941.125
$ws.Range("I16").Value = 2962.75
$ws.Range("J16").Value = 4919.5
$ws.Range("K16").Value = 2962.75
$ws.Range("L16").Value = 4919.5
$ws.Range("M16").Value = -2675.75
$ws.Range("N16").Value = -5493.5
# row 31 (id 44023)
$ws.Range("H31").Value = 4426.95
$ws.Range("I31").Value = 1281.25
$ws.Range("J31").Value = 5213.375
$ws.Range("K31").Value = 1281.25
$ws.Range("L31").Value = 5213.375
$ws.Range("M31").Value = -986.25
$ws.Range("N31").Value = -5803.375
# row 34 (id 44023)
$ws.Range("H34").Value = 4426.95
$ws.Range("I34").Value = 1281.25
$ws.Range("J34").Value = 5213.375
$ws.Range("K34").Value = 1281.25
$ws.Range("L34").Value = 5213.375
$ws.Range("M34").Value = -1079.25
$ws.Range("N34").Value = -5617.375
# row 58 (id 44021)
$ws.Range("H58").Value = 1676.9
$ws.Range("I58").Value = 1284.6364
$ws.Range("K58").Value = 1284.6364
$ws.Range("M58").Value = -1081.6364
# row 62 (id 12580)
$ws.Range("H62").Value = 3741.6667
$ws.Range("I62").Value = 1483.3334
$ws.Range("K62").Value = 1483.3334
$ws.Range("M62").Value = -859.3334
# row 65 (id 12580)
$ws.Range("H65").Value = 3741.6667
$ws.Range("I65").Value = 1483.3334
$ws.Range("K65").Value = 7416.666999999999
$ws.Range("M65").Value = -4296.666999999999
# row 75 (id 11936)
$ws.Range("H75").Value = 30000
$ws.Range("I75").Value = 30000
$ws.Range("K75").Value = 30000
$ws.Range("M75").Value = -29002
# row 76 (id 39765)
$ws.Range("H76").Value = 6233.3335
$ws.Range("I76").Value = 6233.3335
$ws.Range("K76").Value = 6233.3335
$ws.Range("M76").Value = -5918.3335
# row 78 (id 11936)
$ws.Range("H78").Value = 30000
$ws.Range("I78").Value = 30000
$ws.Range("K78").Value = 90000
$ws.Range("M78").Value = -85008
# row 79 (id 39765)
$ws.Range("H79").Value = 6233.3335
$ws.Range("I79").Value = 6233.3335
$ws.Range("K79").Value = 6233.3335
$ws.Range("M79").Value = -5141.3335
# row 86 (id 12584)
$ws.Range("H86").Value = 11474.1
$ws.Range("I86").Value = 10139.357
$ws.Range("K86").Value = 10139.357
$ws.Range("M86").Value = -9016.357
# row 89 (id 12584)
$ws.Range("H89").Value = 11474.1
$ws.Range("I89").Value = 10139.357
$ws.Range("K89").Value = 50696.785
$ws.Range("M89").Value = -45080.785
# row 105 (id 19928)
$ws.Range("H105").Value = 2027.4166
$ws.Range("I105").Value = 2248.3333
$ws.Range("J105").Value = 1806.5
$ws.Range("K105").Value = 2248.3333
$ws.Range("L105").Value = 1806.5
$ws.Range("M105").Value = -501.3332999999998
$ws.Range("N105").Value = -5300.5
# row 113 (id 27691)
$ws.Range("H113").Value = 3941.125
$ws.Range("I113").Value = 2962.75
$ws.Range("J113").Value = 4919.5
$ws.Range("K113").Value = 2962.75
$ws.Range("L113").Value = 4919.5
$ws.Range("M113").Value = -792.75
$ws.Range("N113").Value = -9259.5
# row 132 (id 44019)
$ws.Range("H132").Value = 61457.176
$ws.Range("I132").Value = 64735.75
$ws.Range("K132").Value = 194207.25
$ws.Range("M132").Value = -191677.25
# row 134 (id 44020)
$ws.Range("H134").Value = 2919.7222
$ws.Range("I134").Value = 2125.158
$ws.Range("K134").Value = 6375.474
$ws.Range("M134").Value = -3840.474
# row 136 (id 44021)
$ws.Range("H136").Value = 1676.9
$ws.Range("I136").Value = 1284.6364
$ws.Range("K136").Value = 3853.9092
$ws.Range("M136").Value = -1303.9092
# row 138 (id 42302)
$ws.Range("H138").Value = 97500
$ws.Range("J138").Value = 97500
$ws.Range("L138").Value = 97500
$ws.Range("N138").Value = -107780

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 3 (id 44094)
$ws.Range("H3").Value = 1826
# row 68 (id 12895)
$ws.Range("H68").Value = 2294.7144
$ws.Range("J68").Value = 2804.6667
$ws.Range("L68").Value = 8414.000100000001
$ws.Range("N68").Value = -10036.0001
# row 71 (id 12895)
$ws.Range("H71").Value = 2294.7144
$ws.Range("J71").Value = 2804.6667
$ws.Range("L71").Value = 25242.0003
$ws.Range("N71").Value = -33354.0003
# row 92 (id 19841)
$ws.Range("H92").Value = 550
$ws.Range("J92").Value = 490.375
$ws.Range("L92").Value = 1471.125
$ws.Range("N92").Value = -3967.125
# row 97 (id 19846)
$ws.Range("H97").Value = 212.14285
$ws.Range("I97").Value = 212.14285
$ws.Range("K97").Value = 636.4285500000001
$ws.Range("M97").Value = -140.4285500000001
# row 98 (id 19843)
$ws.Range("H98").Value = 430.8
$ws.Range("I98").Value = 220
$ws.Range("J98").Value = 483.5
$ws.Range("K98").Value = 660
$ws.Range("L98").Value = 1450.5
$ws.Range("M98").Value = 838
$ws.Range("N98").Value = -4446.5
# row 107 (id 27838)
$ws.Range("H107").Value = 1243.75
$ws.Range("J107").Value = 1438.9
$ws.Range("L107").Value = 4316.700000000001
$ws.Range("N107").Value = -8156.700000000001
# row 113 (id 27843)
$ws.Range("H113").Value = 1976.238
$ws.Range("I113").Value = 3247.7273
$ws.Range("J113").Value = 1525.0646
$ws.Range("K113").Value = 9743.1819
$ws.Range("L113").Value = 4575.1938
$ws.Range("M113").Value = -7573.1819
$ws.Range("N113").Value = -8915.193800000001
# row 131 (id 36060)
$ws.Range("H131").Value = 3269.4
$ws.Range("I131").Value = 800.3333
$ws.Range("J131").Value = 6973
$ws.Range("K131").Value = 2400.9999
$ws.Range("L131").Value = 20919
$ws.Range("M131").Value = 2639.0001
$ws.Range("N131").Value = -30999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 69 (id 11891)
$ws.Range("H69").Value = 25000
$ws.Range("J69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("N69").Value = -26498
# row 70 (id 14146)
$ws.Range("H70").Value = 40004600
$ws.Range("I70").Value = 50004000
$ws.Range("K70").Value = 50004000
$ws.Range("M70").Value = -50003730
# row 72 (id 11891)
$ws.Range("H72").Value = 25000
$ws.Range("J72").Value = 25000
$ws.Range("L72").Value = 75000
$ws.Range("N72").Value = -82488
# row 73 (id 14146)
$ws.Range("H73").Value = 40004600
$ws.Range("I73").Value = 50004000
$ws.Range("K73").Value = 50004000
$ws.Range("M73").Value = -50003064
# row 97 (id 19940)
$ws.Range("H97").Value = 994116.7
$ws.Range("I97").Value = 1192521.5
$ws.Range("K97").Value = 1192521.5
$ws.Range("M97").Value = -1192025.5
# row 102 (id 36169)
$ws.Range("H102").Value = 5912519
$ws.Range("I102").Value = 11113540
$ws.Range("J102").Value = 1911733.8
$ws.Range("K102").Value = 11113540
$ws.Range("L102").Value = 1911733.8
$ws.Range("M102").Value = -11111918
$ws.Range("N102").Value = -1914977.8
# row 107 (id 27802)
$ws.Range("H107").Value = 521.7
$ws.Range("I107").Value = 672.53845
$ws.Range("J107").Value = 241.57143
$ws.Range("K107").Value = 672.53845
$ws.Range("L107").Value = 241.57143
$ws.Range("M107").Value = 1247.46155
$ws.Range("N107").Value = -4081.57143
# row 120 (id 26336)
$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -44676
# row 126 (id 36184)
$ws.Range("H126").Value = 5370699.5
$ws.Range("I126").Value = 3035279.2
$ws.Range("K126").Value = 9105837.600000001
$ws.Range("M126").Value = -9103367.600000001
# row 132 (id 44008)
$ws.Range("H132").Value = 2638.9333
$ws.Range("I132").Value = 2230.4736
$ws.Range("K132").Value = 6691.4208
$ws.Range("M132").Value = -4161.4208
# row 133 (id 41854)
$ws.Range("H133").Value = 109199.5
$ws.Range("J133").Value = 109199.5
$ws.Range("L133").Value = 109199.5
$ws.Range("N133").Value = -119319.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7 (id 36249)
$ws.Range("H7").Value = 4142.737
$ws.Range("I7").Value = 2732.125
$ws.Range("J7").Value = 11666
$ws.Range("K7").Value = 2732.125
$ws.Range("L7").Value = 11666
$ws.Range("M7").Value = -2620.125
$ws.Range("N7").Value = -11890
# row 22 (id 5277)
$ws.Range("H22").Value = 89768.60000000001
$ws.Range("I22").Value = 296896
$ws.Range("J22").Value = 999.7143
$ws.Range("K22").Value = 296896
$ws.Range("L22").Value = 999.7143
$ws.Range("M22").Value = -296601
$ws.Range("N22").Value = -1589.7143
# row 27 (id 5277)
$ws.Range("H27").Value = 89768.60000000001
$ws.Range("I27").Value = 296896
$ws.Range("J27").Value = 999.7143
$ws.Range("K27").Value = 296896
$ws.Range("L27").Value = 999.7143
$ws.Range("M27").Value = -296789
$ws.Range("N27").Value = -1213.7143
# row 40 (id 36248)
$ws.Range("H40").Value = 9403.929
$ws.Range("I40").Value = 8200.521000000001
$ws.Range("J40").Value = 14939.6
$ws.Range("K40").Value = 8200.521000000001
$ws.Range("L40").Value = 14939.6
$ws.Range("M40").Value = -8064.521000000001
$ws.Range("N40").Value = -15211.6
# row 54 (id 3177)
$ws.Range("H54").Value = 37520
$ws.Range("I54").Value = 30000
$ws.Range("K54").Value = 30000
$ws.Range("M54").Value = -29356
# row 55 (id 5284)
$ws.Range("H55").Value = 1217.9762
$ws.Range("I55").Value = 896.06665
$ws.Range("J55").Value = 2022.75
$ws.Range("K55").Value = 896.06665
$ws.Range("L55").Value = 2022.75
$ws.Range("M55").Value = -723.06665
$ws.Range("N55").Value = -2368.75
# row 61 (id 27740)
$ws.Range("H61").Value = 22224582
$ws.Range("I61").Value = 22224582
$ws.Range("K61").Value = 22224582
$ws.Range("M61").Value = -22224380
# row 93 (id 19993)
$ws.Range("H93").Value = 12346439
$ws.Range("I93").Value = 15873677
$ws.Range("K93").Value = 15873677
$ws.Range("M93").Value = -15872429
# row 113 (id 27740)
$ws.Range("H113").Value = 22224582
$ws.Range("I113").Value = 22224582
$ws.Range("K113").Value = 22224582
$ws.Range("M113").Value = -22222412
# row 122 (id 36247)
$ws.Range("H122").Value = 5859.769
$ws.Range("I122").Value = 3898.625
$ws.Range("J122").Value = 8997.6
$ws.Range("K122").Value = 11695.875
$ws.Range("L122").Value = 26992.8
$ws.Range("M122").Value = -9245.875
$ws.Range("N122").Value = -31892.8
# row 126 (id 36249)
$ws.Range("H126").Value = 4142.737
$ws.Range("I126").Value = 2732.125
$ws.Range("J126").Value = 11666
$ws.Range("K126").Value = 8196.375
$ws.Range("L126").Value = 34998
$ws.Range("M126").Value = -5726.375
$ws.Range("N126").Value = -39938
# row 132 (id 44058)
$ws.Range("H132").Value = 5575
$ws.Range("I132").Value = 3920
$ws.Range("J132").Value = 8333.333000000001
$ws.Range("K132").Value = 11760
$ws.Range("L132").Value = 24999.999
$ws.Range("M132").Value = -9230
$ws.Range("N132").Value = -30059.999
# row 134 (id 42024)
$ws.Range("H134").Value = 83449.836
$ws.Range("J134").Value = 83449.836
$ws.Range("L134").Value = 83449.836
$ws.Range("N134").Value = -93589.836
# row 136 (id 44060)
$ws.Range("H136").Value = 31301.805
$ws.Range("I136").Value = 36175
$ws.Range("J136").Value = 6935.8335
$ws.Range("K136").Value = 108525
$ws.Range("L136").Value = 20807.5005
$ws.Range("M136").Value = -105975
$ws.Range("N136").Value = -25907.5005
# row 138 (id 42334)
$ws.Range("H138").Value = 71999.664
$ws.Range("J138").Value = 71999.664
$ws.Range("L138").Value = 71999.664
$ws.Range("N138").Value = -82279.664
# row 141 (id 42487)
$ws.Range("H141").Value = 127500
$ws.Range("J141").Value = 127500
$ws.Range("L141").Value = 127500
$ws.Range("N141").Value = -137860

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 17 (id 3539)
$ws.Range("H17").Value = 2500
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null
# row 39 (id 3106)
$ws.Range("H39").Value = 27333.334
$ws.Range("I39").Value = 27333.334
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 27333.334
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -26920.334
$ws.Range("N39").Value = $null
# row 62 (id 12589)
$ws.Range("H62").Value = 26400.2
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248
# row 65 (id 12589)
$ws.Range("H65").Value = 26400.2
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
# row 81 (id 12596)
$ws.Range("H81").Value = 41670216
$ws.Range("I81").Value = 41670216
$ws.Range("K81").Value = 83340432
$ws.Range("M81").Value = -83339371
# row 84 (id 12596)
$ws.Range("H84").Value = 41670216
$ws.Range("I84").Value = 41670216
$ws.Range("K84").Value = 416702160
$ws.Range("M84").Value = -416696856
# row 94 (id 18075)
$ws.Range("H94").Value = 26500
$ws.Range("J94").Value = 26500
$ws.Range("L94").Value = 26500
$ws.Range("N94").Value = -28302
# row 100 (id 19981)
$ws.Range("H100").Value = 675
$ws.Range("I100").Value = 350
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 700
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -159
$ws.Range("N100").Value = -3082
# row 122 (id 36208)
$ws.Range("H122").Value = 2699.5833
$ws.Range("I122").Value = 2542.6562
$ws.Range("K122").Value = 7627.9686
$ws.Range("M122").Value = -5177.9686
# row 126 (id 36210)
$ws.Range("H126").Value = 2436.077
$ws.Range("I126").Value = 2305.3333
$ws.Range("K126").Value = 6915.999899999999
$ws.Range("M126").Value = -4445.999899999999
# row 131 (id 34723)
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080
# row 132 (id 44029)
$ws.Range("H132").Value = 77764616
$ws.Range("I132").Value = 125004200
$ws.Range("J132").Value = 2181277.5
$ws.Range("K132").Value = 375012600
$ws.Range("L132").Value = 6543832.5
$ws.Range("M132").Value = -375010070
$ws.Range("N132").Value = -6548892.5
# row 135 (id 42043)
$ws.Range("H135").Value = 31000
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null
# row 136 (id 44031)
$ws.Range("H136").Value = 3967
$ws.Range("I136").Value = 3996.9744
$ws.Range("J136").Value = 3837.111
$ws.Range("K136").Value = 11990.9232
$ws.Range("L136").Value = 11511.333
$ws.Range("M136").Value = -9440.923200000001
$ws.Range("N136").Value = -16611.333
